$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.934983372688293
$ws.Range("B1").Value = 5.930744647979736
$ws.Range("C1").Value = 3.379580736160278
$ws.Range("D1").Value = 1.441963195800781
$ws.Range("E1").Value = 0.9677280187606812
